$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new rows 17-19 (A index 15-17) have the same column-A number style as existing data rows
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 10: HKL index 8 -> "Gaussian-Quadrature"
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10, 3).Value = 1.013900055707453
$ws.Cells.Item(10, 4).Value = 0.966172179807626
$ws.Cells.Item(10, 5).Value = 1.001089611008914
$ws.Cells.Item(10, 6).Value = 0.9923383171122281
$ws.Cells.Item(10, 7).Value = 1.013900055707453
$ws.Cells.Item(10, 8).Value = 0.966172179807626
$ws.Cells.Item(10, 9).Value = 1.005272623953762
$ws.Cells.Item(10, 10).Value = 0.9889396371195713
$ws.Cells.Item(10, 11).Value = 1.002280705661392
$ws.Cells.Item(10, 12).Value = 0.9785689931768242
$ws.Cells.Item(10, 13).Value = 1.013900055707453
$ws.Cells.Item(10, 14).Value = 0.98363089540827
$ws.Cells.Item(10, 15).Value = 0.9933750409090552
$ws.Cells.Item(10, 16).Value = 0.9935702654434713

# Row 11: HKL index 9 -> "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 3).Value = 1.063951796849523
$ws.Cells.Item(11, 4).Value = 0.8876575370010542
$ws.Cells.Item(11, 5).Value = 1.027876986452726
$ws.Cells.Item(11, 6).Value = 0.9714994836172389
$ws.Cells.Item(11, 7).Value = 1.063951796849523
$ws.Cells.Item(11, 8).Value = 0.8876575370010542
$ws.Cells.Item(11, 9).Value = 1.036140927605159
$ws.Cells.Item(11, 10).Value = 0.974632232627624
$ws.Cells.Item(11, 11).Value = 1.02
$ws.Cells.Item(11, 12).Value = 0.93
$ws.Cells.Item(11, 13).Value = 1.063951796849523
$ws.Cells.Item(11, 14).Value = 0.9577672617268902
$ws.Cells.Item(11, 15).Value = 0.9877464509801355
$ws.Cells.Item(11, 16).Value = 0.9889698705191656

# Row 12: HKL index 10 -> "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 3).Value = 1.063958518950846
$ws.Cells.Item(12, 4).Value = 0.8876903629624887
$ws.Cells.Item(12, 5).Value = 1.027876299235434
$ws.Cells.Item(12, 6).Value = 0.9714815803602485
$ws.Cells.Item(12, 7).Value = 1.063958518950846
$ws.Cells.Item(12, 8).Value = 0.8876903629624887
$ws.Cells.Item(12, 9).Value = 1.036135016889764
$ws.Cells.Item(12, 10).Value = 0.9746476094286329
$ws.Cells.Item(12, 11).Value = 1.02
$ws.Cells.Item(12, 12).Value = 0.93
$ws.Cells.Item(12, 13).Value = 1.063958518950846
$ws.Cells.Item(12, 14).Value = 0.9577833310989614
$ws.Cells.Item(12, 15).Value = 0.9877516903772543
$ws.Cells.Item(12, 16).Value = 0.9889736734784267

# Row 13: HKL index 11 -> "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 3).Value = 1.063921991618633
$ws.Cells.Item(13, 4).Value = 0.8876765645429671
$ws.Cells.Item(13, 5).Value = 1.027883038813431
$ws.Cells.Item(13, 6).Value = 0.9715006841308189
$ws.Cells.Item(13, 7).Value = 1.063921991618633
$ws.Cells.Item(13, 8).Value = 0.8876765645429671
$ws.Cells.Item(13, 9).Value = 1.036143487531012
$ws.Cells.Item(13, 10).Value = 0.9746383766646206
$ws.Cells.Item(13, 11).Value = 1.02
$ws.Cells.Item(13, 12).Value = 0.93
$ws.Cells.Item(13, 13).Value = 1.063921991618633
$ws.Cells.Item(13, 14).Value = 0.9577798016781991
$ws.Cells.Item(13, 15).Value = 0.9877455697764624
$ws.Cells.Item(13, 16).Value = 0.9889705179126853

# Row 14: HKL index 12 -> "NoRotation-tilt60deg"
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14, 3).Value = 1.077723999999999
$ws.Cells.Item(14, 4).Value = 0.8893799999999999
$ws.Cells.Item(14, 5).Value = 1.024576000000002
$ws.Cells.Item(14, 6).Value = 0.9708320000000004
$ws.Cells.Item(14, 7).Value = 1.077723999999999
$ws.Cells.Item(14, 8).Value = 0.8893799999999999
$ws.Cells.Item(14, 9).Value = 1.032912
$ws.Cells.Item(14, 10).Value = 0.9729119999999997
$ws.Cells.Item(14, 11).Value = 1.02
$ws.Cells.Item(14, 12).Value = 0.93
$ws.Cells.Item(14, 13).Value = 1.077723999999999
$ws.Cells.Item(14, 14).Value = 0.9569780000000008
$ws.Cells.Item(14, 15).Value = 0.9906280000000005
$ws.Cells.Item(14, 16).Value = 0.9897920000000001

# Row 15: HKL index 13 -> "Rotation-NoTilt"
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15, 3).Value = 1.09
$ws.Cells.Item(15, 4).Value = 0.8943125
$ws.Cells.Item(15, 5).Value = 1.02
$ws.Cells.Item(15, 6).Value = 0.97
$ws.Cells.Item(15, 7).Value = 1.09
$ws.Cells.Item(15, 8).Value = 0.8943125
$ws.Cells.Item(15, 9).Value = 1.03
$ws.Cells.Item(15, 10).Value = 0.9691749999999989
$ws.Cells.Item(15, 11).Value = 1.018887499999998
$ws.Cells.Item(15, 12).Value = 0.93
$ws.Cells.Item(15, 13).Value = 1.09
$ws.Cells.Item(15, 14).Value = 0.9571562499999999
$ws.Cells.Item(15, 15).Value = 0.993578125
$ws.Cells.Item(15, 16).Value = 0.9902968749999996

# Row 16: HKL index 14 -> "Rotation-60detTilt"
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16, 3).Value = 1.049025518591993
$ws.Cells.Item(16, 4).Value = 0.9369967702016022
$ws.Cells.Item(16, 5).Value = 1.010238069964799
$ws.Cells.Item(16, 6).Value = 0.9809546158080009
$ws.Cells.Item(16, 7).Value = 1.049025518591993
$ws.Cells.Item(16, 8).Value = 0.9369967702016022
$ws.Cells.Item(16, 9).Value = 1.015746453299195
$ws.Cells.Item(16, 10).Value = 0.9808642073600049
$ws.Cells.Item(16, 11).Value = 1.009644355481595
$ws.Cells.Item(16, 12).Value = 0.9579551954944003
$ws.Cells.Item(16, 13).Value = 1.049025518591993
$ws.Cells.Item(16, 14).Value = 0.9736174200832004
$ws.Cells.Item(16, 15).Value = 0.9943037436415987
$ws.Cells.Item(16, 16).Value = 0.9926781482751987

# Row 17: HKL index 15 -> "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 3).Value = 0.9959893203102888
$ws.Cells.Item(17, 4).Value = 0.9960232981127799
$ws.Cells.Item(17, 5).Value = 0.9961152031849769
$ws.Cells.Item(17, 6).Value = 0.9958622372371321
$ws.Cells.Item(17, 7).Value = 0.9959893203102888
$ws.Cells.Item(17, 8).Value = 0.9960232981127799
$ws.Cells.Item(17, 9).Value = 0.9959035428369406
$ws.Cells.Item(17, 10).Value = 0.9961190647439339
$ws.Cells.Item(17, 11).Value = 0.9965041205508318
$ws.Cells.Item(17, 12).Value = 0.9959122778873029
$ws.Cells.Item(17, 13).Value = 0.9959893203102888
$ws.Cells.Item(17, 14).Value = 0.9960692506488784
$ws.Cells.Item(17, 15).Value = 0.9959975147112944
$ws.Cells.Item(17, 16).Value = 0.9960536331080234

# Row 18: HKL index 16 -> "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 3).Value = 0.9955447106344336
$ws.Cells.Item(18, 4).Value = 0.9963770346229226
$ws.Cells.Item(18, 5).Value = 0.9952506645184048
$ws.Cells.Item(18, 6).Value = 0.9959054672325155
$ws.Cells.Item(18, 7).Value = 0.9955447106344336
$ws.Cells.Item(18, 8).Value = 0.9963770346229226
$ws.Cells.Item(18, 9).Value = 0.9949280483052149
$ws.Cells.Item(18, 10).Value = 0.9972030977905006
$ws.Cells.Item(18, 11).Value = 0.9952478369001835
$ws.Cells.Item(18, 12).Value = 0.9989733663909722
$ws.Cells.Item(18, 13).Value = 0.9955447106344336
$ws.Cells.Item(18, 14).Value = 0.9958138495706637
$ws.Cells.Item(18, 15).Value = 0.9957694692520691
$ws.Cells.Item(18, 16).Value = 0.9961787782993935

# Row 19: HKL index 17 -> "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 3).Value = 0.9879164157093315
$ws.Cells.Item(19, 4).Value = 1.007756344560311
$ws.Cells.Item(19, 5).Value = 0.9938769154117232
$ws.Cells.Item(19, 6).Value = 0.9983587992948711
$ws.Cells.Item(19, 7).Value = 0.9879164157093315
$ws.Cells.Item(19, 8).Value = 1.007756344560311
$ws.Cells.Item(19, 9).Value = 0.9918861618717072
$ws.Cells.Item(19, 10).Value = 0.9987427621236625
$ws.Cells.Item(19, 11).Value = 0.9934149041475179
$ws.Cells.Item(19, 12).Value = 1.003705397843503
$ws.Cells.Item(19, 13).Value = 0.9879164157093315
$ws.Cells.Item(19, 14).Value = 1.000816629986017
$ws.Cells.Item(19, 15).Value = 0.9969771187440591
$ws.Cells.Item(19, 16).Value = 0.9969572126203283
